$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark from the tail of the document;
# it will be re-created at the end of the new last paragraph below.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Position right after the text of the last paragraph ("Format koji
# cuvam ..."), before its paragraph mark, and split off a new blank
# paragraph.
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

# This new trailing blank paragraph becomes the separator before the
# new heading.
$blankPara1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $blankPara1.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# The next paragraph holds the new "Graph Classes" heading.
$headingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$hr = $headingPara.Range
$hr.Collapse(0)
$hr.InsertAfter("Graph Classes")
$headingPara.Style = "Heading 1"

$hr2 = $headingPara.Range
$hr2.Collapse(0)
$hr2.InsertParagraphAfter()

# The body paragraph discussing GraphClasses vs RandomGraphClasses.
$bodyPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bodyPara.Style = "Normal"
$br = $bodyPara.Range
$br.Collapse(0)
$br.InsertAfter("Vjerovatno bi bilo ljepse da GraphClasses budu klase za koje imamo klasifikatore koje racunamo, a da izdvojimo ove slozenije za generisanje random grafova kao RandomGraphClasses, jer oni nisu bas klase, nego su cesto kombinacija klasa ili neka klasa koja se generise zadavanjem odredjenih parametara. Drugi nacin")
$br.Collapse(0)
$br.InsertAfter(" je da ih flagujemo ne")
$br.Collapse(0)
$br.InsertAfter("kako samo da su za kalkulisanje, odnosno za generisanje (sto vec postoji).")
$br.Collapse(0)
$br.InsertAfter(" Mozda mozemo za sad sa ovim drugim pristupom.")

# --- Re-create the "_GoBack" bookmark right after the body paragraph's
# text, before its paragraph mark. Bookmarks.Add is unreliable when the
# collapsed target position is exactly "paragraph end - 1" (i.e. right
# before a paragraph mark), so we temporarily append a one-character
# placeholder, add the bookmark just before that placeholder (a safe
# position), then remove the placeholder again.
$phRange = $bodyPara.Range
$phRange.Collapse(0)
[void]$phRange.MoveEnd(1, -1)
$phRange.Collapse(0)
$phRange.InsertAfter("X")

$targetPos = $bodyPara.Range.End - 2
$bmRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$placeholderRange = $d.Range($bodyPara.Range.End - 2, $bodyPara.Range.End - 1)
$placeholderRange.Delete()

# Trailing blank paragraph after the new body paragraph.
$br2 = $bodyPara.Range
$br2.Collapse(0)
$br2.InsertParagraphAfter()
$trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailingPara.Style = "Normal"
